$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("F2").Value = 11845.35
$ws.Range("AG2").Value = 71495.85000000001

# Row 3 - Bibi Cell Ponta Negra
$ws.Range("E3").Value = 5592
$ws.Range("F3").Value = 3002
$ws.Range("AG3").Value = 16812.52

# Row 4 - Bibi Cell Vieiralves
$ws.Range("F4").Value = 4238
$ws.Range("AG4").Value = 15529

# Row 5 - Bibi Cell Manauara
$ws.Range("E5").Value = 1819
$ws.Range("F5").Value = 2503
$ws.Range("AG5").Value = 12115

# Row 6 - total
$ws.Range("E6").Value = 36189.46
$ws.Range("F6").Value = 21588.35
$ws.Range("AG6").Value = 115952.37
